$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new header row at the top of Sheet2, pushing the existing
# data (and the shared-string indices that go with it) down by one row.
$ws.Rows.Item(1).Insert()

# Populate the new header row with the "menu" / "menu_class" labels.
$ws.Range("A1").Value = "menu"
$ws.Range("B1").Value = "menu_class"

# Reflect the final cursor position recorded in the saved file.
$ws.Range("D9").Select()
